$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("F1").Value = "DATAEVENTO"

# New numeric (date-serial) values for column F, rows 2-7
$ws.Cells.Item(2, 6).Value = 45489.70232528935
$ws.Cells.Item(3, 6).Value = 45489.70143650463
$ws.Cells.Item(4, 6).Value = 45489.70181196759
$ws.Cells.Item(5, 6).Value = 45489.702215509256
$ws.Cells.Item(6, 6).Value = 45489.702560578706
$ws.Cells.Item(7, 6).Value = 45489.70259893518
